# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Column E (Periodo Mora) for rows 16-27 previously ran 2401..2412 ascending.
# The data was rebuilt so the periods now run in descending order (2412..2401).
$periods = @(2412, 2411, 2410, 2409, 2408, 2407, 2406, 2405, 2404, 2403, 2402, 2401)
$row = 16
foreach ($p in $periods) {
    $ws.Cells.Item($row, 5).Value = "$p"
    $row++
}

# Column F (Valor Mora): the 19067 value that belonged to period 2412 (row 27)
# now travels with that period to row 16; the remaining rows keep 52000, so
# only the two previously-differing cells swap.
$ws.Range("F16").Value = 19067
$ws.Range("F27").Value = 52000
